$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update K column (최종점수 / final score) values
$ws.Range("K2").Value = 57.5
$ws.Range("K3").Value = 55.7
$ws.Range("K4").Value = 54.5
$ws.Range("K5").Value = 54.5

# Update N column (MACRO_SCORE) values
$ws.Range("N2").Value = 51.53902399942638
$ws.Range("N3").Value = 51.53902399942638
$ws.Range("N4").Value = 51.53902399942638
$ws.Range("N5").Value = 51.53902399942638
